# "Code changes for white slot"
#
# Adds a new TestData test case (TC2) that is a copy of the existing TC1
# row, while TC1's NegotiatedBy/Network are updated to "Business Affairs "
# / "Telemundo". Also updates the active sheet/selection bookkeeping:
# TestData becomes the active sheet (selection D3) and Windows keeps
# selection D4 but is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")
$wsWindows  = $wb.Worksheets.Item("Windows")

# Duplicate the existing TC1 row (row 2) down into a new row 3, carrying
# over its values/formatting so the new TC2 row starts out identical to
# the original TC1 row.
$wsTestData.Rows.Item(2).Copy()
$wsTestData.Rows.Item(3).Insert()

# The new row 3 is the TC2 test case; give it its own TcNo.
$wsTestData.Range("A3").Value() = "TC2"

# The original row (now row 2, TC1) gets a new negotiator / network.
$wsTestData.Range("D2").Value() = "Business Affairs "
$wsTestData.Range("E2").Value() = "Telemundo"

# Update the Windows sheet selection (it stops being the active tab, but
# keeps a remembered selection of D4).
$wsWindows.Activate()
$wsWindows.Range("D4").Select()

# TestData becomes the active sheet, with D3 (the new row) selected.
$wsTestData.Activate()
$wsTestData.Range("D3").Select()
